$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 10 ("fossil_routes"), shifting rows 10-24 down to 11-25.
$ws.Rows.Item(10).Insert()

# Fill in the new row 10 with the new parameter "chemical_recycling_pyrolysis" = TRUE
$ws.Cells.Item(10, 1).Value = "chemical_recycling_pyrolysis"
$ws.Cells.Item(10, 2).Value = $true

# Update the sheet dimension/used range is handled automatically by Excel.
